$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# D-column cells are forced to Text format so numeric-looking strings
# (e.g. thousand-separated prices like '63.925.84') are preserved exactly.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.925.84"
$ws.Range("E2").Value = "  -1.74%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.144.28"
$ws.Range("E3").Value = "  -1.24%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.98"
$ws.Range("E5").Value = "  -2.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.91"
$ws.Range("E6").Value = "  -4.10%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.139.47"
$ws.Range("E8").Value = "  -1.34%  "

$ws.Range("E9").Value = "  -0.27%  "

$ws.Range("E10").Value = "  -3.15%  "

$ws.Range("E11").Value = "  -2.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000253"
$ws.Range("E13").Value = "  -3.94%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.87"
$ws.Range("E14").Value = "  -3.42%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.665.81"
$ws.Range("E15").Value = "  -1.08%  "

$ws.Range("E16").Value = "  +2.53%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.948.76"
$ws.Range("E17").Value = "  -1.68%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.142.59"
$ws.Range("E18").Value = "  -1.28%  "

$ws.Range("E19").Value = "  -2.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "487.27"
$ws.Range("E20").Value = "  +0.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.68"
$ws.Range("E21").Value = "  -0.80%  "

$ws.Range("E22").Value = "  -2.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.74"
$ws.Range("E23").Value = "  -3.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.33"
$ws.Range("E24").Value = "  +4.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.22"
$ws.Range("E25").Value = "  -5.46%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("E27").Value = "  -2.65%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.18"
$ws.Range("E28").Value = "  -6.31%  "

$ws.Range("E29").Value = "  -2.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.06"
$ws.Range("E30").Value = "  -3.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.44"
$ws.Range("E31").Value = "  +2.21%  "

$ws.Range("E32").Value = "  -7.94%  "

$ws.Range("E33").Value = "  -0.08%  "

$ws.Range("E34").Value = "  -3.50%  "

$ws.Range("E35").Value = "  -3.28%  "

$ws.Range("E36").Value = "  -0.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.65"
$ws.Range("E37").Value = "  -1.16%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0743"
$ws.Range("E38").Value = "  -6.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.92"
$ws.Range("E39").Value = "  -8.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0397"
$ws.Range("E40").Value = "  -1.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "431.15"
$ws.Range("E41").Value = "  -8.40%  "

$ws.Range("E42").Value = "  -0.51%  "

$ws.Range("E43").Value = "  -0.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.912.32"
$ws.Range("E44").Value = "  +1.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.259"
$ws.Range("E45").Value = "  -4.40%  "

$ws.Range("E46").Value = "  -7.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.39"
$ws.Range("E47").Value = "  -3.38%  "

$ws.Range("E48").Value = "  -0.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.81"
$ws.Range("E51").Value = "  -0.47%  "

# Rows 49 and 50: Stellar and InjectiveProtocol swapped positions
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.115"
$ws.Range("E49").Value = "  -0.37%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.70"
$ws.Range("E50").Value = "  -4.76%  "
